$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6; existing rows 6..58 shift down to 7..59.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Cells.Item(6, 1).Value = 9
$ws.Cells.Item(6, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(6, 3).Value = "Metropolitana"
$ws.Cells.Item(6, 4).Value = Get-Date -Year 2023 -Month 5 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(6, 5).Value = 13
$ws.Cells.Item(6, 6).Value = 100112010
$ws.Cells.Item(6, 7).Value = "Achicoria"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 90
$ws.Cells.Item(6, 11).Value = 7000
$ws.Cells.Item(6, 12).Value = 7000
$ws.Cells.Item(6, 13).Value = 7000
$ws.Cells.Item(6, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(6, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(6, 16).Value = 438
$ws.Cells.Item(6, 17).Value = 16
$ws.Cells.Item(6, 18).Value = "Hortaliza"

# Match the date-cell formatting/style used by the rest of column D.
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
